# QA and update files to console
# The "IP provider" row label is renamed to "ISP Line" and highlighted in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch")

# C3 used to read "IP provider" (shared string referencing the Chinese/English
# localization pair). Replace it with the corrected "ISP Line" label and mark
# it in red so reviewers notice the QA fix. (Name/Size already match the
# sheet's base Arial 12pt font, so only the color needs to change.)
$cell = $ws.Range("C3")
$cell.Value = "ISP Line"
$cell.Font.Color = 255

# Switch the sheet to Portrait orientation (matches the saved page setup).
$ws.PageSetup.Orientation = 1

# Leave the saved cursor position on C10, matching the workbook's last view.
[void]$ws.Range("C10").Select()
